$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 5917
$ws.Range("G2").Value = 138
$ws.Range("G3").Value = 168
$ws.Range("F4").Value = 1113
$ws.Range("F5").Value = 1057
$ws.Range("F8").Value = 52
$ws.Range("F10").Value = 65
$ws.Range("F11").Value = 31
$ws.Range("F12").Value = 28
$ws.Range("F13").Value = 2078
$ws.Range("F14").Value = 1527
$ws.Range("F15").Value = 1145
$ws.Range("F18").Value = 443
$ws.Range("F19").Value = 668
$ws.Range("F20").Value = 236
$ws.Range("F21").Value = 1074
$ws.Range("F24").Value = 3772
$ws.Range("F28").Value = 171
$ws.Range("F29").Value = 54
$ws.Range("F30").Value = 541
$ws.Range("F32").Value = 55
$ws.Range("F33").Value = 25
$ws.Range("F36").Value = 862
$ws.Range("F37").Value = 111
$ws.Range("F39").Value = 90
$ws.Range("F40").Value = 93

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = "已停售"

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G2").Value = "已停售"
$ws.Range("F3").Value = 5917
$ws.Range("G3").Value = 138
$ws.Range("G4").Value = 168
$ws.Range("F5").Value = 1113
$ws.Range("F8").Value = 1057
$ws.Range("F13").Value = 52
$ws.Range("F15").Value = 65
$ws.Range("F16").Value = 31
$ws.Range("F18").Value = 28
$ws.Range("F19").Value = 2078
$ws.Range("F20").Value = 1527
$ws.Range("F21").Value = 1145
$ws.Range("F24").Value = 443
$ws.Range("F26").Value = 668
$ws.Range("F27").Value = 236
$ws.Range("F28").Value = 1074
$ws.Range("F30").Value = 3772
$ws.Range("F34").Value = 171
$ws.Range("F35").Value = 54
$ws.Range("F36").Value = 541
$ws.Range("F38").Value = 55
$ws.Range("F39").Value = 25
$ws.Range("F42").Value = 862
$ws.Range("F43").Value = 111
$ws.Range("F45").Value = 90
$ws.Range("F46").Value = 93
